$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.138179
$ws.Range("H2").Value = 0.414537
$ws.Range("I2").Value = 0.03953416978071796
$ws.Range("J2").Value = 0.03953416978071796
$ws.Range("M2").Value = 0.9949089999999999
$ws.Range("N2").Value = 2.984727
$ws.Range("O2").Value = 0.1476822527339178
$ws.Range("P2").Value = 0.1476822527339178
$ws.Range("Q2").Value = 0.137475530711
$ws.Range("R2").Value = 1.237279776399
$ws.Range("S2").Value = 0.005838495253181605
$ws.Range("T2").Value = 0.005838495253181605
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.138179
$ws.Range("H3").Value = 0.414537
$ws.Range("I3").Value = 0.03953416978071796
$ws.Range("J3").Value = 0.03953416978071796
$ws.Range("O3").Value = 0.2453919293791607
$ws.Range("P3").Value = 0.2453919293791607
$ws.Range("Q3").Value = 0.2284322259383333
$ws.Range("R3").Value = 2.055890033445
$ws.Range("S3").Value = 0.009701366198893691
$ws.Range("T3").Value = 0.009701366198893691
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.138179
$ws.Range("H4").Value = 0.414537
$ws.Range("I4").Value = 0.03953416978071796
$ws.Range("J4").Value = 0.03953416978071796
$ws.Range("M4").Value = 1.748891
$ws.Range("N4").Value = 5.246673
$ws.Range("O4").Value = 0.2596017954064887
$ws.Range("P4").Value = 0.2596017954064887
$ws.Range("Q4").Value = 0.241660009489
$ws.Range("R4").Value = 2.174940085401
$ws.Range("S4").Value = 0.01026314145497933
$ws.Range("T4").Value = 0.01026314145497933
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.138179
$ws.Range("H5").Value = 0.414537
$ws.Range("I5").Value = 0.03953416978071796
$ws.Range("J5").Value = 0.03953416978071796
$ws.Range("M5").Value = 0.7268083333333334
$ws.Range("N5").Value = 2.180425
$ws.Range("O5").Value = 0.1078859392893731
$ws.Range("P5").Value = 0.1078859392893731
$ws.Range("Q5").Value = 0.1004296486916667
$ws.Range("R5").Value = 0.903866838225
$ws.Range("S5").Value = 0.004265181040818306
$ws.Range("T5").Value = 0.004265181040818306
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.138179
$ws.Range("H6").Value = 0.414537
$ws.Range("I6").Value = 0.03953416978071796
$ws.Range("J6").Value = 0.03953416978071796
$ws.Range("M6").Value = 1.613051666666667
$ws.Range("N6").Value = 4.839155
$ws.Range("O6").Value = 0.2394380831910597
$ws.Range("P6").Value = 0.2394380831910597
$ws.Range("Q6").Value = 0.2228898662483333
$ws.Range("R6").Value = 2.006008796235
$ws.Range("S6").Value = 0.009465985832845024
$ws.Range("T6").Value = 0.009465985832845024
$ws.Range("I7").Value = 0.9514265220751211
$ws.Range("J7").Value = 0.9514265220751212
$ws.Range("M7").Value = 0.9949089999999999
$ws.Range("N7").Value = 2.984727
$ws.Range("O7").Value = 0.1476822527339178
$ws.Range("P7").Value = 0.1476822527339178
$ws.Range("Q7").Value = 3.308476358054
$ws.Range("R7").Value = 29.776287222486
$ws.Range("S7").Value = 0.1405088120908504
$ws.Range("T7").Value = 0.1405088120908505
$ws.Range("I8").Value = 0.9514265220751211
$ws.Range("J8").Value = 0.9514265220751212
$ws.Range("O8").Value = 0.2453919293791607
$ws.Range("P8").Value = 0.2453919293791607
$ws.Range("S8").Value = 0.2334723899145186
$ws.Range("T8").Value = 0.2334723899145187
$ws.Range("I9").Value = 0.9514265220751211
$ws.Range("J9").Value = 0.9514265220751212
$ws.Range("M9").Value = 1.748891
$ws.Range("N9").Value = 5.246673
$ws.Range("O9").Value = 0.2596017954064887
$ws.Range("P9").Value = 0.2596017954064887
$ws.Range("Q9").Value = 5.815772624746
$ws.Range("R9").Value = 52.341953622714
$ws.Range("S9").Value = 0.2469920333280527
$ws.Range("T9").Value = 0.2469920333280527
$ws.Range("I10").Value = 0.9514265220751211
$ws.Range("J10").Value = 0.9514265220751212
$ws.Range("M10").Value = 0.7268083333333334
$ws.Range("N10").Value = 2.180425
$ws.Range("O10").Value = 0.1078859392893731
$ws.Range("P10").Value = 0.1078859392893731
$ws.Range("Q10").Value = 2.416932792516667
$ws.Range("R10").Value = 21.75239513265
$ws.Range("S10").Value = 0.1026455439988959
$ws.Range("T10").Value = 0.1026455439988959
$ws.Range("I11").Value = 0.9514265220751211
$ws.Range("J11").Value = 0.9514265220751212
$ws.Range("M11").Value = 1.613051666666667
$ws.Range("N11").Value = 4.839155
$ws.Range("O11").Value = 0.2394380831910597
$ws.Range("P11").Value = 0.2394380831910597
$ws.Range("Q11").Value = 5.364051690643333
$ws.Range("R11").Value = 48.27646521579
$ws.Range("S11").Value = 0.2278077427428034
$ws.Range("T11").Value = 0.2278077427428035
$ws.Range("G12").Value = 0.031594
$ws.Range("H12").Value = 0.09478200000000001
$ws.Range("I12").Value = 0.009039308144160858
$ws.Range("J12").Value = 0.009039308144160858
$ws.Range("M12").Value = 0.9949089999999999
$ws.Range("N12").Value = 2.984727
$ws.Range("O12").Value = 0.1476822527339178
$ws.Range("P12").Value = 0.1476822527339178
$ws.Range("Q12").Value = 0.031433154946
$ws.Range("R12").Value = 0.282898394514
$ws.Range("S12").Value = 0.001334945389885725
$ws.Range("T12").Value = 0.001334945389885725
$ws.Range("G13").Value = 0.031594
$ws.Range("H13").Value = 0.09478200000000001
$ws.Range("I13").Value = 0.009039308144160858
$ws.Range("J13").Value = 0.009039308144160858
$ws.Range("O13").Value = 0.2453919293791607
$ws.Range("P13").Value = 0.2453919293791607
$ws.Range("Q13").Value = 0.05222998969666667
$ws.Range("R13").Value = 0.47006990727
$ws.Range("S13").Value = 0.002218173265748394
$ws.Range("T13").Value = 0.002218173265748394
$ws.Range("G14").Value = 0.031594
$ws.Range("H14").Value = 0.09478200000000001
$ws.Range("I14").Value = 0.009039308144160858
$ws.Range("J14").Value = 0.009039308144160858
$ws.Range("M14").Value = 1.748891
$ws.Range("N14").Value = 5.246673
$ws.Range("O14").Value = 0.2596017954064887
$ws.Range("P14").Value = 0.2596017954064887
$ws.Range("Q14").Value = 0.05525446225400001
$ws.Range("R14").Value = 0.497290160286
$ws.Range("S14").Value = 0.002346620623456654
$ws.Range("T14").Value = 0.002346620623456654
$ws.Range("G15").Value = 0.031594
$ws.Range("H15").Value = 0.09478200000000001
$ws.Range("I15").Value = 0.009039308144160858
$ws.Range("J15").Value = 0.009039308144160858
$ws.Range("M15").Value = 0.7268083333333334
$ws.Range("N15").Value = 2.180425
$ws.Range("O15").Value = 0.1078859392893731
$ws.Range("P15").Value = 0.1078859392893731
$ws.Range("Q15").Value = 0.02296278248333334
$ws.Range("R15").Value = 0.20666504235
$ws.Range("S15").Value = 0.0009752142496588743
$ws.Range("T15").Value = 0.0009752142496588743
$ws.Range("G16").Value = 0.031594
$ws.Range("H16").Value = 0.09478200000000001
$ws.Range("I16").Value = 0.009039308144160858
$ws.Range("J16").Value = 0.009039308144160858
$ws.Range("M16").Value = 1.613051666666667
$ws.Range("N16").Value = 4.839155
$ws.Range("O16").Value = 0.2394380831910597
$ws.Range("P16").Value = 0.2394380831910597
$ws.Range("Q16").Value = 0.05096275435666667
$ws.Range("R16").Value = 0.45866478921
$ws.Range("S16").Value = 0.002164354615411211
$ws.Range("T16").Value = 0.002164354615411211
